$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "success" header in D1, matching the style of the other headers ---
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value = "success"

# --- Per-row success flag ("1" success / "0" failure) for rows 2..35 ---
$success = @{
    2  = "1"
    3  = "1"
    4  = "1"
    5  = "1"
    6  = "1"
    7  = "0"
    8  = "1"
    9  = "1"
    10 = "0"
    11 = "1"
    12 = "1"
    13 = "1"
    14 = "0"
    15 = "1"
    16 = "1"
    17 = "0"
    18 = "0"
    19 = "0"
    20 = "0"
    21 = "1"
    22 = "1"
    23 = "1"
    24 = "1"
    25 = "1"
    26 = "1"
    27 = "0"
    28 = "1"
    29 = "1"
    30 = "1"
    31 = "1"
    32 = "0"
    33 = "1"
    34 = "1"
    35 = "0"
}

for ($row = 2; $row -le 35; $row++) {
    $cell = $ws.Cells.Item($row, 4)
    $val = $success[$row]
    # Build the literal text value through a formula so Excel stores it
    # as shared-string text ("1"/"0") instead of inferring a number, then
    # bake the formula result down to a plain value via paste-special.
    $cell.Formula = "=""" + $val + """"
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

$excel.CutCopyMode = 0

Write-Host "done"
